$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.426.32"
$ws.Range("D3").Value = "1.553.73"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'210.45"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "'0.482"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'24.12"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "'0.242"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "'0.0891"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "1.775.88"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "1.553.10"
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "28.428.65"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "'3.62"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "'61.03"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "'228.61"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("D20").Value = "0.0₃0672"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'3.89"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").Value = "'151.26"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "'14.72"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D29").Value = "'6.22"
$ws.Range("E29").Value = "  -3.08%  "
$ws.Range("D30").Value = "'0.0466"
$ws.Range("E30").Value = "  -3.11%  "
$ws.Range("E31").Value = "  -4.61%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").Value = "1.384.01"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("D36").Value = "'1.47"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "'0.772"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'5.35"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").Value = "'61.75"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").Value = "1.688.51"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E48").Value = "  -7.53%  "
$ws.Range("D49").Value = "'85.04"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("D50").Value = "'42.73"
$ws.Range("E50").Value = "  +6.90%  "
$ws.Range("E51").Value = "  -2.35%  "
